# Generate Report for Handback
# Adds a new handed-back file (7b2d61c3-682e-4462-bb61-f80f93f3d795.md) as
# row 4 to the "Overview", "zh-cn" and "de-de" worksheets, wires up the new
# hyperlinks and grows the tables / dimensions to match.

$wb = $excel.ActiveWorkbook

$guid = "7b2d61c3-682e-4462-bb61-f80f93f3d795"
$fileName = "$guid.md"
$pathAndName = "e2e\$guid.md"
$displayNameOnly = "$guid.md"

$zhXlf = "$guid.681b89667a631ab356e7a1d05e3abee8b421c308.zh-cn.xlf"
$deXlf = "$guid.681b89667a631ab356e7a1d05e3abee8b421c308.de-de.xlf"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = $fileName
$wsOverview.Range("B4").Value = $pathAndName
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("E4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G4").Value = "2016-09-04 02:48:08"
$wsOverview.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7b2d61c3682e4462bb61f80f93f3d795682e446/e2e/$guid.md", [Type]::Missing, [Type]::Missing, $pathAndName) | Out-Null

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G4"))

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A4").Value = $fileName
$wsZh.Range("B4").Value = ".md"
$wsZh.Range("C4").Value = "Handed back: in sync with en-US"
$wsZh.Range("D4").Value = "e2e"
$wsZh.Range("E4").Value = "ht"
$wsZh.Range("F4").Value = "True"
$wsZh.Range("G4").Value = $zhXlf
$wsZh.Range("H4").Value = "2016-09-04 02:47:59"
$wsZh.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("I4").Value = $fileName
$wsZh.Range("J4").Value = $zhXlf
$wsZh.Range("K4").Value = "2016-09-04 02:48:33"
$wsZh.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsZh.Range("L4").Value = "'"
$wsZh.Range("L4").Style = "Normal"
$wsZh.Range("M4").Value = "True"
$wsZh.Range("N4").Value = "'"
$wsZh.Range("N4").Style = "Normal"
$wsZh.Range("O4").Value = "False"
$wsZh.Range("P4").Value = "'"
$wsZh.Range("P4").Style = "Normal"

$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7b2d61c3682e4462bb61f80f93f3d795682e446/e2e/$guid.md", [Type]::Missing, [Type]::Missing, $displayNameOnly) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/681b89667a631ab356e7a1d05e3abee8b421c308/e2e/$guid.md", [Type]::Missing, [Type]::Missing, $displayNameOnly) | Out-Null

$loZh = $wsZh.ListObjects.Item(1)
$loZh.Resize($wsZh.Range("A1:P4"))

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A4").Value = $fileName
$wsDe.Range("B4").Value = ".md"
$wsDe.Range("C4").Value = "Handed back: in sync with en-US"
$wsDe.Range("D4").Value = "e2e"
$wsDe.Range("E4").Value = "ht"
$wsDe.Range("F4").Value = "True"
$wsDe.Range("G4").Value = $deXlf
$wsDe.Range("H4").Value = "2016-09-04 02:48:08"
$wsDe.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("I4").Value = $fileName
$wsDe.Range("J4").Value = $deXlf
$wsDe.Range("K4").Value = "2016-09-04 02:48:41"
$wsDe.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsDe.Range("L4").Value = "'"
$wsDe.Range("L4").Style = "Normal"
$wsDe.Range("M4").Value = "True"
$wsDe.Range("N4").Value = "'"
$wsDe.Range("N4").Style = "Normal"
$wsDe.Range("O4").Value = "False"
$wsDe.Range("P4").Value = "'"
$wsDe.Range("P4").Style = "Normal"

$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7b2d61c3682e4462bb61f80f93f3d795682e446/e2e/$guid.md", [Type]::Missing, [Type]::Missing, $displayNameOnly) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/681b89667a631ab356e7a1d05e3abee8b421c308/e2e/$guid.md", [Type]::Missing, [Type]::Missing, $displayNameOnly) | Out-Null

$loDe = $wsDe.ListObjects.Item(1)
$loDe.Resize($wsDe.Range("A1:P4"))
